$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column B ("Year of Treatment"); this shifts C:Q left to B:P
$ws.Range("B:B").Delete()

# Append ".global" suffix to all header cells in row 1 except A1 ("Country")
for ($col = 2; $col -le 16; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = [string]$cell.Value2 + ".global"
}
